$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 20919.2
$ws.Range("I96").Value = 25949
$ws.Range("K96").Value = 77847
$ws.Range("M96").Value = -76474
$ws.Range("H106").Value = 30269.428
$ws.Range("I106").Value = 31147.666
$ws.Range("K106").Value = 31147.666
$ws.Range("M106").Value = -30516.666
$ws.Range("H113").Value = 4400
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 4400
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").Value = 4400
$ws.Range("N113").Value = -10908
$ws.Range("H137").Value = 1827.0834
$ws.Range("I137").Value = 1886.2222
$ws.Range("J137").Value = 1649.6666
$ws.Range("K137").Value = 5658.6666
$ws.Range("L137").Value = 4948.9998
$ws.Range("M137").Value = -3108.6666
$ws.Range("N137").Value = -10048.9998
$ws.Range("H138").Value = 5620.522
$ws.Range("J138").Value = 5909.263
$ws.Range("L138").Value = 17727.789
$ws.Range("N138").Value = -28007.789

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1586
$ws.Range("I4").Value = 344
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 344
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = -228
$ws.Range("N4").Value = -2232
$ws.Range("H32").Value = 2096.6304
$ws.Range("I32").Value = 2098.4666
$ws.Range("K32").Value = 2098.4666
$ws.Range("M32").Value = -1811.4666
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").ClearContents()
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = 0
$ws.Range("H55").Value = 50000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 50000
$ws.Range("K55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("M55").Value = 50000
$ws.Range("N55").Value = -50630
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("N123").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1436
$ws.Range("I20").Value = 1501.9
$ws.Range("J20").Value = 1216.3334
$ws.Range("K20").Value = 1501.9
$ws.Range("L20").Value = 1216.3334
$ws.Range("M20").Value = -1254.9
$ws.Range("N20").Value = -1710.3334
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8783.333000000001
$ws.Range("I4").Value = 8800
$ws.Range("J4").Value = 8775
$ws.Range("K4").Value = 8800
$ws.Range("L4").Value = 8775
$ws.Range("M4").Value = -8688
$ws.Range("N4").Value = -8999
$ws.Range("H58").Value = 1861.3334
$ws.Range("I58").Value = 1536
$ws.Range("K58").Value = 1536
$ws.Range("M58").Value = -1333
$ws.Range("H62").Value = 4409.9
$ws.Range("I62").Value = 4299.8335
$ws.Range("K62").Value = 4299.8335
$ws.Range("M62").Value = -3675.8335
$ws.Range("H65").Value = 4409.9
$ws.Range("I65").Value = 4299.8335
$ws.Range("K65").Value = 21499.1675
$ws.Range("M65").Value = -18379.1675
$ws.Range("H99").Value = 3250
$ws.Range("I99").Value = 1500
$ws.Range("K99").Value = 1500
$ws.Range("M99").Value = -2
$ws.Range("H126").Value = 3250
$ws.Range("I126").Value = 1500
$ws.Range("K126").Value = 4500
$ws.Range("M126").Value = -2030
$ws.Range("H132").Value = 3115.25
$ws.Range("I132").Value = 3230.1428
$ws.Range("K132").Value = 9690.428400000001
$ws.Range("M132").Value = -7160.428400000001
$ws.Range("H136").Value = 1861.3334
$ws.Range("I136").Value = 1536
$ws.Range("K136").Value = 4608
$ws.Range("M136").Value = -2058

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 442223.8
$ws.Range("I2").Value = 366706
$ws.Range("K2").Value = 2200236
$ws.Range("M2").Value = -2200123
$ws.Range("H4").Value = 18833432
$ws.Range("I4").Value = 18833432
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 56500296
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -56500184
$ws.Range("H7").Value = 277.14285
$ws.Range("I7").Value = 152.33333
$ws.Range("J7").Value = 370.75
$ws.Range("K7").Value = 456.99999
$ws.Range("L7").Value = 1112.25
$ws.Range("M7").Value = -344.99999
$ws.Range("N7").Value = -1336.25
$ws.Range("H132").Value = 920.9167
$ws.Range("I132").Value = 552
$ws.Range("K132").Value = 4968
$ws.Range("M132").Value = -2438

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 19500
$ws.Range("I5").Value = 10000
$ws.Range("K5").Value = 10000
$ws.Range("M5").Value = -9888
$ws.Range("H10").Value = 9250
$ws.Range("J10").Value = 9250
$ws.Range("L10").Value = 9250
$ws.Range("N10").Value = -9588
$ws.Range("H102").Value = 4489.7144
$ws.Range("I102").Value = 4312.909
$ws.Range("K102").Value = 4312.909
$ws.Range("M102").Value = -2690.909
$ws.Range("H126").Value = 4500
$ws.Range("I126").Value = 4500
$ws.Range("K126").Value = 13500
$ws.Range("M126").Value = -11030
$ws.Range("H140").Value = 143848.5
$ws.Range("J140").Value = 143848.5
$ws.Range("L140").Value = 143848.5
$ws.Range("N140").Value = -154208.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 139966.33
$ws.Range("J2").Value = 19899
$ws.Range("L2").Value = 19899
$ws.Range("N2").Value = -20123
$ws.Range("H61").Value = 5500.8
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H100").Value = 1958.35
$ws.Range("I100").Value = 1786.7646
$ws.Range("K100").Value = 1786.7646
$ws.Range("M100").Value = -1245.7646
$ws.Range("H113").Value = 5500.8
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H132").Value = 4869.9414
$ws.Range("I132").Value = 5384.5386
$ws.Range("J132").Value = 3197.5
$ws.Range("K132").Value = 16153.6158
$ws.Range("L132").Value = 9592.5
$ws.Range("M132").Value = -13623.6158
$ws.Range("N132").Value = -14652.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 150000
$ws.Range("J2").Value = 150000
$ws.Range("L2").Value = 150000
$ws.Range("N2").Value = -150224
$ws.Range("H81").Value = 1000864
$ws.Range("I81").Value = 1034
$ws.Range("K81").Value = 2068
$ws.Range("M81").Value = -1007
$ws.Range("H84").Value = 1000864
$ws.Range("I84").Value = 1034
$ws.Range("K84").Value = 10340
$ws.Range("M84").Value = -5036
$ws.Range("H126").Value = 2156.75
$ws.Range("J126").Value = 1050
$ws.Range("L126").Value = 3150
$ws.Range("N126").Value = -8090
